$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr/text cells) so exact
# formatted strings (trailing zeros, percent signs) are preserved verbatim.
$cells = @{
    "D2" = '294.52'
    "E2" = '1.70%'
    "E3" = '0.34%'
    "D4" = '4.934'
    "E4" = '1.38%'
    "D5" = '0.07353'
    "E5" = '3.11%'
    "D6" = '2.314'
    "E6" = '28.03%'
    "D7" = '7.695'
    "E7" = '0.60%'
    "D8" = '3.760'
    "E8" = '-0.50%'
    "D9" = '0.9131'
    "E9" = '1.88%'
    "D10" = '0.1686'
    "E10" = '2.78%'
    "D11" = '0.08183'
    "E11" = '8.50%'
    "D12" = '0.08277'
    "E12" = '2.96%'
    "D13" = '0.03129'
    "E13" = '4.74%'
    "D14" = '0.1009'
    "E14" = '0.91%'
    "D15" = '0.001511'
    "E15" = '0.42%'
    "D16" = '0.005740'
    "E16" = '-2.00%'
    "E17" = '0.68%'
    "D18" = '2.074'
    "E18" = '-1.57%'
    "D19" = '0.3329'
    "E19" = '1.61%'
    "E20" = '0.46%'
    "D21" = '3.972'
    "E21" = '-6.98%'
    "E22" = '4.94%'
    "D23" = '0.04551'
    "E23" = '1.32%'
    "D24" = '0.001211'
    "E24" = '-0.14%'
    "D25" = '0.004338'
    "E25" = '-6.89%'
    "E26" = '3.88%'
    "D27" = '0.0003395'
    "E39" = '-1.81%'
    "D40" = '0.04448'
    "E40" = '2.21%'
    "D41" = '0.007360'
    "E41" = '-0.65%'
    "D42" = '0.008837'
    "D43" = '0.1323'
    "E43" = '1.23%'
    "D44" = '0.002081'
    "E44" = '3.78%'
    "D45" = '0.009106'
    "E45" = '-11.97%'
    "D46" = '0.00005904'
    "E46" = '0.99%'
    "E47" = '-0.13%'
    "E48" = '1.08%'
    "E50" = '-0.13%'
    "E51" = '-0.13%'
}

foreach ($ref in $cells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$ref]
}
